$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2211538461538461
$ws.Range("C2").Value = 0.4711538461538461
$ws.Range("J2").Value = 0.03525641025641026
$ws.Range("P2").Value = 0.1762820512820513
$ws.Range("S2").Value = 0.09615384615384616
$ws.Range("B3").Value = 0.006578947368421052
$ws.Range("C3").Value = 0.0131578947368421
$ws.Range("J3").Value = 0.05263157894736842
$ws.Range("P3").Value = 0.7039473684210527
$ws.Range("S3").Value = 0.2236842105263158
$ws.Range("B6").Value = 0.06324110671936758
$ws.Range("D6").Value = 0.003952569169960474
$ws.Range("F6").Value = 0.05928853754940711
$ws.Range("J6").Value = 0.2490118577075099
$ws.Range("O6").Value = 0.007905138339920948
$ws.Range("Q6").Value = 0.1067193675889328
$ws.Range("R6").Value = 0.07509881422924901
$ws.Range("S6").Value = 0.4347826086956522
$ws.Range("B7").Value = 0.07575757575757576
$ws.Range("D7").Value = 0.0202020202020202
$ws.Range("E7").Value = 0.005050505050505051
$ws.Range("F7").Value = 0.08585858585858586
$ws.Range("J7").Value = 0.1060606060606061
$ws.Range("O7").Value = 0.01515151515151515
$ws.Range("Q7").Value = 0.1363636363636364
$ws.Range("R7").Value = 0.0707070707070707
$ws.Range("S7").Value = 0.4848484848484849
$ws.Range("B8").Value = 0.1014492753623188
$ws.Range("D8").Value = 0.01449275362318841
$ws.Range("F8").Value = 0.06340579710144928
$ws.Range("J8").Value = 0.1177536231884058
$ws.Range("O8").Value = 0.03442028985507246
$ws.Range("Q8").Value = 0.1902173913043478
$ws.Range("R8").Value = 0.05797101449275362
$ws.Range("S8").Value = 0.4202898550724637
$ws.Range("B9").Value = 0.08962264150943396
$ws.Range("D9").Value = 0.01415094339622642
$ws.Range("F9").Value = 0.0660377358490566
$ws.Range("J9").Value = 0.09433962264150944
$ws.Range("O9").Value = 0.0330188679245283
$ws.Range("Q9").Value = 0.1745283018867924
$ws.Range("R9").Value = 0.0660377358490566
$ws.Range("S9").Value = 0.4622641509433962
$ws.Range("B10").Value = 0.09759759759759759
$ws.Range("D10").Value = 0.01501501501501501
$ws.Range("F10").Value = 0.08033033033033032
$ws.Range("J10").Value = 0.1388888888888889
$ws.Range("O10").Value = 0.01726726726726727
$ws.Range("Q10").Value = 0.1876876876876877
$ws.Range("R10").Value = 0.05855855855855856
$ws.Range("S10").Value = 0.4046546546546547
$ws.Range("G11").Value = 0.124113475177305
$ws.Range("J11").Value = 0.09219858156028368
$ws.Range("K11").Value = 0.198581560283688
$ws.Range("L11").Value = 0.5780141843971631
$ws.Range("S11").Value = 0.007092198581560284
$ws.Range("G12").Value = 0.7797619047619048
$ws.Range("J12").Value = 0.1488095238095238
$ws.Range("K12").Value = 0.0119047619047619
$ws.Range("L12").Value = 0.02380952380952381
$ws.Range("S12").Value = 0.03571428571428571
$ws.Range("G13").Value = 0.75
$ws.Range("J13").Value = 0.2307692307692308
$ws.Range("S13").Value = 0.01923076923076923
$ws.Range("F15").Value = 0.02702702702702703
$ws.Range("H15").Value = 0.1486486486486487
$ws.Range("I15").Value = 0.1036036036036036
$ws.Range("J15").Value = 0.3108108108108108
$ws.Range("K15").Value = 0.04054054054054054
$ws.Range("M15").Value = 0.009009009009009009
$ws.Range("O15").Value = 0.06306306306306306
$ws.Range("S15").Value = 0.2972972972972973
$ws.Range("F16").Value = 0.00558659217877095
$ws.Range("H16").Value = 0.1843575418994413
$ws.Range("I16").Value = 0.07262569832402235
$ws.Range("J16").Value = 0.4357541899441341
$ws.Range("K16").Value = 0.1173184357541899
$ws.Range("M16").Value = 0.00558659217877095
$ws.Range("N16").Value = 0.00558659217877095
$ws.Range("O16").Value = 0.03910614525139665
$ws.Range("S16").Value = 0.1340782122905028
$ws.Range("F17").Value = 0.025
$ws.Range("H17").Value = 0.2340909090909091
$ws.Range("I17").Value = 0.09772727272727273
$ws.Range("J17").Value = 0.3909090909090909
$ws.Range("K17").Value = 0.08863636363636364
$ws.Range("M17").Value = 0.01590909090909091
$ws.Range("O17").Value = 0.04318181818181818
$ws.Range("S17").Value = 0.1045454545454545
$ws.Range("F18").Value = 0.006369426751592357
$ws.Range("H18").Value = 0.1847133757961783
$ws.Range("I18").Value = 0.1210191082802548
$ws.Range("J18").Value = 0.3949044585987261
$ws.Range("K18").Value = 0.1019108280254777
$ws.Range("M18").Value = 0.01273885350318471
$ws.Range("O18").Value = 0.03821656050955414
$ws.Range("S18").Value = 0.1401273885350318
$ws.Range("F19").Value = 0.0107095046854083
$ws.Range("H19").Value = 0.2396251673360107
$ws.Range("I19").Value = 0.07697456492637215
$ws.Range("J19").Value = 0.35809906291834
$ws.Range("K19").Value = 0.09236947791164658
$ws.Range("M19").Value = 0.02811244979919679
$ws.Range("N19").Value = 0.001338688085676037
$ws.Range("O19").Value = 0.06425702811244979
$ws.Range("S19").Value = 0.1285140562248996
